# Update profit sheet with the new row for 09/11/2025
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 24
$newRow = $lastRow + 1

# Write the date as text (matching the existing style of the other date
# cells in column A, which store the date as a literal string rather than
# a date serial number). The leading apostrophe forces Excel to treat the
# value as text instead of auto-converting it to a date.
$ws.Range("A" + $newRow).Value = "'09/11/2025"

# Match the (lack of) explicit styling used by the surrounding date cells.
$ws.Range("A" + $newRow).Style = $ws.Range("A" + $lastRow).Style

# Write the profit figure as a plain number.
$ws.Range("B" + $newRow).Value = 15253.51
